$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.415243666666667
$ws.Range("H2").Value = 7.245730999999999
$ws.Range("I2").Value = 0.006221624451646611
$ws.Range("J2").Value = 0.006221624451646612
$ws.Range("M2").Value = 28.31444233333334
$ws.Range("N2").Value = 84.94332700000001
$ws.Range("O2").Value = 0.2747173016130739
$ws.Range("P2").Value = 0.2747173016130739
$ws.Range("Q2").Value = 68.38627752078189
$ws.Range("R2").Value = 615.4764976870371
$ws.Range("S2").Value = 0.001709187881006278
$ws.Range("T2").Value = 0.001709187881006278
$ws.Range("G3").Value = 2.415243666666667
$ws.Range("H3").Value = 7.245730999999999
$ws.Range("I3").Value = 0.006221624451646611
$ws.Range("J3").Value = 0.006221624451646612
$ws.Range("O3").Value = 0.2090339131726295
$ws.Range("P3").Value = 0.2090339131726295
$ws.Range("Q3").Value = 52.03549653968411
$ws.Range("R3").Value = 468.3194688571569
$ws.Range("S3").Value = 0.001300530505418206
$ws.Range("T3").Value = 0.001300530505418206
$ws.Range("G4").Value = 2.415243666666667
$ws.Range("H4").Value = 7.245730999999999
$ws.Range("I4").Value = 0.006221624451646611
$ws.Range("J4").Value = 0.006221624451646612
$ws.Range("M4").Value = 5.413469333333334
$ws.Range("N4").Value = 16.240408
$ws.Range("O4").Value = 0.0525235026743817
$ws.Range("P4").Value = 0.0525235026743817
$ws.Range("Q4").Value = 13.07484752202756
$ws.Range("R4").Value = 117.673627698248
$ws.Range("S4").Value = 0.0003267815085250594
$ws.Range("T4").Value = 0.0003267815085250594
$ws.Range("G5").Value = 2.415243666666667
$ws.Range("H5").Value = 7.245730999999999
$ws.Range("I5").Value = 0.006221624451646611
$ws.Range("J5").Value = 0.006221624451646612
$ws.Range("M5").Value = 47.79503400000001
$ws.Range("N5").Value = 143.385102
$ws.Range("O5").Value = 0.4637252825399149
$ws.Range("P5").Value = 0.4637252825399149
$ws.Range("Q5").Value = 115.436653166618
$ws.Range("R5").Value = 1038.929878499562
$ws.Range("S5").Value = 0.002885124556697068
$ws.Range("T5").Value = 0.002885124556697068
$ws.Range("I6").Value = 0.8933025543886721
$ws.Range("J6").Value = 0.8933025543886722
$ws.Range("M6").Value = 28.31444233333334
$ws.Range("N6").Value = 84.94332700000001
$ws.Range("O6").Value = 0.2747173016130739
$ws.Range("P6").Value = 0.2747173016130739
$ws.Range("Q6").Value = 9818.920584041221
$ws.Range("R6").Value = 88370.285256371
$ws.Range("S6").Value = 0.2454056672657222
$ws.Range("T6").Value = 0.2454056672657223
$ws.Range("I7").Value = 0.8933025543886721
$ws.Range("J7").Value = 0.8933025543886722
$ws.Range("O7").Value = 0.2090339131726295
$ws.Range("P7").Value = 0.2090339131726295
$ws.Range("S7").Value = 0.1867305285909698
$ws.Range("T7").Value = 0.1867305285909699
$ws.Range("I8").Value = 0.8933025543886721
$ws.Range("J8").Value = 0.8933025543886722
$ws.Range("M8").Value = 5.413469333333334
$ws.Range("N8").Value = 16.240408
$ws.Range("O8").Value = 0.0525235026743817
$ws.Range("P8").Value = 0.0525235026743817
$ws.Range("Q8").Value = 1877.290212619383
$ws.Range("R8").Value = 16895.61191357444
$ws.Range("S8").Value = 0.04691937910446543
$ws.Range("T8").Value = 0.04691937910446543
$ws.Range("I9").Value = 0.8933025543886721
$ws.Range("J9").Value = 0.8933025543886722
$ws.Range("M9").Value = 47.79503400000001
$ws.Range("N9").Value = 143.385102
$ws.Range("O9").Value = 0.4637252825399149
$ws.Range("P9").Value = 0.4637252825399149
$ws.Range("Q9").Value = 16574.42649347429
$ws.Range("R9").Value = 149169.8384412686
$ws.Range("S9").Value = 0.4142469794275146
$ws.Range("T9").Value = 0.4142469794275147
$ws.Range("G10").Value = 38.75388733333333
$ws.Range("H10").Value = 116.261662
$ws.Range("I10").Value = 0.09982932006284441
$ws.Range("J10").Value = 0.09982932006284444
$ws.Range("M10").Value = 28.31444233333334
$ws.Range("N10").Value = 84.94332700000001
$ws.Range("O10").Value = 0.2747173016130739
$ws.Range("P10").Value = 0.2747173016130739
$ws.Range("Q10").Value = 1097.294708092164
$ws.Range("R10").Value = 9875.652372829476
$ws.Range("S10").Value = 0.02742484142953252
$ws.Range("T10").Value = 0.02742484142953253
$ws.Range("G11").Value = 38.75388733333333
$ws.Range("H11").Value = 116.261662
$ws.Range("I11").Value = 0.09982932006284441
$ws.Range("J11").Value = 0.09982932006284444
$ws.Range("O11").Value = 0.2090339131726295
$ws.Range("P11").Value = 0.2090339131726295
$ws.Range("Q11").Value = 834.9376081859683
$ws.Range("R11").Value = 7514.438473673715
$ws.Range("S11").Value = 0.02086771342209926
$ws.Range("T11").Value = 0.02086771342209926
$ws.Range("G12").Value = 38.75388733333333
$ws.Range("H12").Value = 116.261662
$ws.Range("I12").Value = 0.09982932006284441
$ws.Range("J12").Value = 0.09982932006284444
$ws.Range("M12").Value = 5.413469333333334
$ws.Range("N12").Value = 16.240408
$ws.Range("O12").Value = 0.0525235026743817
$ws.Range("P12").Value = 0.0525235026743817
$ws.Range("Q12").Value = 209.7929806264551
$ws.Range("R12").Value = 1888.136825638096
$ws.Range("S12").Value = 0.005243385559302515
$ws.Range("T12").Value = 0.005243385559302517
$ws.Range("G13").Value = 38.75388733333333
$ws.Range("H13").Value = 116.261662
$ws.Range("I13").Value = 0.09982932006284441
$ws.Range("J13").Value = 0.09982932006284444
$ws.Range("M13").Value = 47.79503400000001
$ws.Range("N13").Value = 143.385102
$ws.Range("O13").Value = 0.4637252825399149
$ws.Range("P13").Value = 0.4637252825399149
$ws.Range("Q13").Value = 1852.243362728836
$ws.Range("R13").Value = 16670.19026455953
$ws.Range("S13").Value = 0.04629337965191011
$ws.Range("T13").Value = 0.04629337965191013
$ws.Range("G14").Value = 0.2509726666666667
$ws.Range("H14").Value = 0.752918
$ws.Range("I14").Value = 0.0006465010968368635
$ws.Range("J14").Value = 0.0006465010968368636
$ws.Range("M14").Value = 28.31444233333334
$ws.Range("N14").Value = 84.94332700000001
$ws.Range("O14").Value = 0.2747173016130739
$ws.Range("P14").Value = 0.2747173016130739
$ws.Range("Q14").Value = 7.106151097576223
$ws.Range("R14").Value = 63.95535987818601
$ws.Range("S14").Value = 0.0001776050368129158
$ws.Range("T14").Value = 0.0001776050368129158
$ws.Range("G15").Value = 0.2509726666666667
$ws.Range("H15").Value = 0.752918
$ws.Range("I15").Value = 0.0006465010968368635
$ws.Range("J15").Value = 0.0006465010968368636
$ws.Range("O15").Value = 0.2090339131726295
$ws.Range("P15").Value = 0.2090339131726295
$ws.Range("Q15").Value = 5.407109646171778
$ws.Range("R15").Value = 48.663986815546
$ws.Range("S15").Value = 0.0001351406541422067
$ws.Range("T15").Value = 0.0001351406541422067
$ws.Range("G16").Value = 0.2509726666666667
$ws.Range("H16").Value = 0.752918
$ws.Range("I16").Value = 0.0006465010968368635
$ws.Range("J16").Value = 0.0006465010968368636
$ws.Range("M16").Value = 5.413469333333334
$ws.Range("N16").Value = 16.240408
$ws.Range("O16").Value = 0.0525235026743817
$ws.Range("P16").Value = 0.0525235026743817
$ws.Range("Q16").Value = 1.358632834504889
$ws.Range("R16").Value = 12.227695510544
$ws.Range("S16").Value = 0.0000339565020887017
$ws.Range("T16").Value = 0.00003395650208870171
$ws.Range("G17").Value = 0.2509726666666667
$ws.Range("H17").Value = 0.752918
$ws.Range("I17").Value = 0.0006465010968368635
$ws.Range("J17").Value = 0.0006465010968368636
$ws.Range("M17").Value = 47.79503400000001
$ws.Range("N17").Value = 143.385102
$ws.Range("O17").Value = 0.4637252825399149
$ws.Range("P17").Value = 0.4637252825399149
$ws.Range("Q17").Value = 1852.243362728836
$ws.Range("R17").Value = 16670.19026455953
$ws.Range("S17").Value = 0.04629337965191011
$ws.Range("T17").Value = 0.04629337965191013
